# The "Delivery" programs sheet used the word "Mixed" to describe a delivery
# method that is both in-person and online. The author decided to rename
# this option to "Both" for clarity, so every occurrence of "Mixed" in the
# "Services" worksheet must become "Both" (commit: "Changed mixed to both
# in combined programs").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Services")

# Select the "Delivery" column (P) — this mirrors how the edit was made by
# hand (selecting the whole column before running Find & Replace) and is
# reflected in the saved selection state of the sheet.
$colP = $ws.Columns.Item(16)
$colP.Select()

# Replace every "Mixed" value with "Both" within that selection/column.
$colP.Replace("Mixed", "Both", 1, 1, $false, $false, $false)
